# Auto-generated edit script: updates market-price derived cells
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# per the scheduled-runner data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 4881.6665  # H38
$ws.Cells.Item(38, 10).Value = 12000  # J38
$ws.Cells.Item(38, 12).Value = 36000  # L38
$ws.Cells.Item(38, 14).Value = -36744  # N38
$ws.Cells.Item(63, 8).Value = 47500  # H63
$ws.Cells.Item(63, 10).Value = 47500  # J63
$ws.Cells.Item(63, 12).Value = 47500  # L63
$ws.Cells.Item(63, 14).Value = -48748  # N63
$ws.Cells.Item(66, 8).Value = 47500  # H66
$ws.Cells.Item(66, 10).Value = 47500  # J66
$ws.Cells.Item(66, 12).Value = 142500  # L66
$ws.Cells.Item(66, 14).Value = -148740  # N66
$ws.Cells.Item(98, 8).Value = 50000350  # H98
$ws.Cells.Item(98, 9).Value = 55555884  # I98
$ws.Cells.Item(98, 11).Value = 55555884  # K98
$ws.Cells.Item(98, 13).Value = -55554386  # M98
$ws.Cells.Item(105, 8).Value = 80685.09  # H105
$ws.Cells.Item(105, 10).Value = 80685.09  # J105
$ws.Cells.Item(105, 12).Value = 80685.09  # L105
$ws.Cells.Item(105, 14).Value = -87673.09  # N105
$ws.Cells.Item(106, 8).Value = 15311.111  # H106
$ws.Cells.Item(106, 9).Value = 2599.6667  # I106
$ws.Cells.Item(106, 10).Value = 21666.834  # J106
$ws.Cells.Item(106, 11).Value = 2599.6667  # K106
$ws.Cells.Item(106, 12).Value = 21666.834  # L106
$ws.Cells.Item(106, 13).Value = -1968.6667  # M106
$ws.Cells.Item(106, 14).Value = -22928.834  # N106
$ws.Cells.Item(122, 8).Value = 50000350  # H122
$ws.Cells.Item(122, 9).Value = 55555884  # I122
$ws.Cells.Item(122, 11).Value = 166667652  # K122
$ws.Cells.Item(122, 13).Value = -166665202  # M122
$ws.Cells.Item(131, 8).Value = 3644.611  # H131
$ws.Cells.Item(131, 9).Value = 1467  # I131
$ws.Cells.Item(131, 10).Value = 7999.8335  # J131
$ws.Cells.Item(131, 11).Value = 4401  # K131
$ws.Cells.Item(131, 12).Value = 23999.5005  # L131
$ws.Cells.Item(131, 13).Value = 639  # M131
$ws.Cells.Item(131, 14).Value = -34079.50049999999  # N131
$ws.Cells.Item(132, 8).Value = 1924.7693  # H132
$ws.Cells.Item(132, 9).Value = 1556.0426  # I132
$ws.Cells.Item(132, 11).Value = 4668.1278  # K132
$ws.Cells.Item(132, 13).Value = -2138.1278  # M132
$ws.Cells.Item(138, 8).Value = 3217.7795  # H138
$ws.Cells.Item(138, 9).Value = 2054.9  # I138
$ws.Cells.Item(138, 10).Value = 3814.1282  # J138
$ws.Cells.Item(138, 11).Value = 6164.700000000001  # K138
$ws.Cells.Item(138, 12).Value = 11442.3846  # L138
$ws.Cells.Item(138, 13).Value = -1024.700000000001  # M138
$ws.Cells.Item(138, 14).Value = -21722.3846  # N138
$ws.Cells.Item(141, 8).Value = 4860.857  # H141
$ws.Cells.Item(141, 9).Value = 5019.4614  # I141
$ws.Cells.Item(141, 11).Value = 15058.3842  # K141
$ws.Cells.Item(141, 13).Value = -9878.3842  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7700752.5  # H32
$ws.Cells.Item(32, 9).Value = 8337176.5  # I32
$ws.Cells.Item(32, 11).Value = 8337176.5  # K32
$ws.Cells.Item(32, 13).Value = -8336889.5  # M32
$ws.Cells.Item(45, 8).Value = 2461.6  # H45
$ws.Cells.Item(45, 9).Value = 2311.5  # I45
$ws.Cells.Item(45, 11).Value = 2311.5  # K45
$ws.Cells.Item(45, 13).Value = -1934.5  # M45
$ws.Cells.Item(132, 8).Value = 5280.1763  # H132
$ws.Cells.Item(132, 9).Value = 2177.5134  # I132
$ws.Cells.Item(132, 10).Value = 13480.071  # J132
$ws.Cells.Item(132, 11).Value = 6532.540199999999  # K132
$ws.Cells.Item(132, 12).Value = 40440.213  # L132
$ws.Cells.Item(132, 13).Value = -4002.540199999999  # M132
$ws.Cells.Item(132, 14).Value = -45500.213  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3201.25  # H107
$ws.Cells.Item(107, 9).Value = 3285  # I107
$ws.Cells.Item(107, 11).Value = 3285  # K107
$ws.Cells.Item(107, 13).Value = -1365  # M107
$ws.Cells.Item(113, 8).Value = 8888  # H113
$ws.Cells.Item(113, 9).Value = 8888  # I113
$ws.Cells.Item(113, 11).Value = 8888  # K113
$ws.Cells.Item(113, 13).Value = -6718  # M113

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 19833.334  # H16
$ws.Cells.Item(16, 9).Value = 20800  # I16
$ws.Cells.Item(16, 10).Value = 15000  # J16
$ws.Cells.Item(16, 11).Value = 20800  # K16
$ws.Cells.Item(16, 12).Value = 15000  # L16
$ws.Cells.Item(16, 13).Value = -20513  # M16
$ws.Cells.Item(16, 14).Value = -15574  # N16
$ws.Cells.Item(31, 8).Value = 590701.8  # H31
$ws.Cells.Item(31, 9).Value = 10180.889  # I31
$ws.Cells.Item(31, 11).Value = 10180.889  # K31
$ws.Cells.Item(31, 13).Value = -9885.888999999999  # M31
$ws.Cells.Item(34, 8).Value = 590701.8  # H34
$ws.Cells.Item(34, 9).Value = 10180.889  # I34
$ws.Cells.Item(34, 11).Value = 10180.889  # K34
$ws.Cells.Item(34, 13).Value = -9978.888999999999  # M34
$ws.Cells.Item(58, 8).Value = 1650.6  # H58
$ws.Cells.Item(58, 9).Value = 1650.6  # I58
$ws.Cells.Item(58, 11).Value = 1650.6  # K58
$ws.Cells.Item(58, 13).Value = -1447.6  # M58
$ws.Cells.Item(63, 8).Value = 49090.332  # H63
$ws.Cells.Item(63, 10).Value = 49090.332  # J63
$ws.Cells.Item(63, 12).Value = 49090.332  # L63
$ws.Cells.Item(63, 14).Value = -50462.332  # N63
$ws.Cells.Item(64, 8).Value = 34666.668  # H64
$ws.Cells.Item(64, 10).Value = 34666.668  # J64
$ws.Cells.Item(64, 12).Value = 34666.668  # L64
$ws.Cells.Item(64, 14).Value = -35162.668  # N64
$ws.Cells.Item(66, 8).Value = 49090.332  # H66
$ws.Cells.Item(66, 10).Value = 49090.332  # J66
$ws.Cells.Item(66, 12).Value = 147270.996  # L66
$ws.Cells.Item(66, 14).Value = -154134.996  # N66
$ws.Cells.Item(67, 8).Value = 34666.668  # H67
$ws.Cells.Item(67, 10).Value = 34666.668  # J67
$ws.Cells.Item(67, 12).Value = 34666.668  # L67
$ws.Cells.Item(67, 14).Value = -36382.668  # N67
$ws.Cells.Item(113, 8).Value = 19833.334  # H113
$ws.Cells.Item(113, 9).Value = 20800  # I113
$ws.Cells.Item(113, 10).Value = 15000  # J113
$ws.Cells.Item(113, 11).Value = 20800  # K113
$ws.Cells.Item(113, 12).Value = 15000  # L113
$ws.Cells.Item(113, 13).Value = -18630  # M113
$ws.Cells.Item(113, 14).Value = -19340  # N113
$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 12).Value = 0  # L124
$ws.Cells.Item(124, 14).ClearContents()  # N124
$ws.Cells.Item(134, 8).Value = 457261.78  # H134
$ws.Cells.Item(134, 9).Value = 528040  # I134
$ws.Cells.Item(134, 11).Value = 1584120  # K134
$ws.Cells.Item(134, 13).Value = -1581585  # M134
$ws.Cells.Item(136, 8).Value = 1650.6  # H136
$ws.Cells.Item(136, 9).Value = 1650.6  # I136
$ws.Cells.Item(136, 11).Value = 4951.799999999999  # K136
$ws.Cells.Item(136, 13).Value = -2401.799999999999  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 30376.25  # H51
$ws.Cells.Item(51, 9).Value = 10000  # I51
$ws.Cells.Item(51, 10).Value = 37168.332  # J51
$ws.Cells.Item(51, 11).Value = 30000  # K51
$ws.Cells.Item(51, 12).Value = 111504.996  # L51
$ws.Cells.Item(51, 13).Value = -29540  # M51
$ws.Cells.Item(51, 14).Value = -112424.996  # N51
$ws.Cells.Item(68, 8).Value = 38936.332  # H68
$ws.Cells.Item(68, 10).Value = 45360.086  # J68
$ws.Cells.Item(68, 12).Value = 136080.258  # L68
$ws.Cells.Item(68, 14).Value = -137702.258  # N68
$ws.Cells.Item(71, 8).Value = 38936.332  # H71
$ws.Cells.Item(71, 10).Value = 45360.086  # J71
$ws.Cells.Item(71, 12).Value = 408240.774  # L71
$ws.Cells.Item(71, 14).Value = -416352.774  # N71
$ws.Cells.Item(76, 8).Value = 4800  # H76
$ws.Cells.Item(76, 9).Value = 3000  # I76
$ws.Cells.Item(76, 11).Value = 9000  # K76
$ws.Cells.Item(76, 13).Value = -8617  # M76
$ws.Cells.Item(79, 8).Value = 4800  # H79
$ws.Cells.Item(79, 9).Value = 3000  # I79
$ws.Cells.Item(79, 11).Value = 9000  # K79
$ws.Cells.Item(79, 13).Value = -7674  # M79

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1900  # H80
$ws.Cells.Item(80, 9).Value = 2500  # I80
$ws.Cells.Item(80, 10).Value = 1300  # J80
$ws.Cells.Item(80, 11).Value = 2500  # K80
$ws.Cells.Item(80, 12).Value = 1300  # L80
$ws.Cells.Item(80, 13).Value = -1502  # M80
$ws.Cells.Item(80, 14).Value = -3296  # N80
$ws.Cells.Item(83, 8).Value = 1900  # H83
$ws.Cells.Item(83, 9).Value = 2500  # I83
$ws.Cells.Item(83, 10).Value = 1300  # J83
$ws.Cells.Item(83, 11).Value = 12500  # K83
$ws.Cells.Item(83, 12).Value = 6500  # L83
$ws.Cells.Item(83, 13).Value = -7508  # M83
$ws.Cells.Item(83, 14).Value = -16484  # N83
$ws.Cells.Item(97, 8).Value = 1413.1578  # H97
$ws.Cells.Item(97, 9).Value = 1575.1428  # I97
$ws.Cells.Item(97, 11).Value = 1575.1428  # K97
$ws.Cells.Item(97, 13).Value = -1079.1428  # M97
$ws.Cells.Item(122, 8).Value = 2083.6785  # H122
$ws.Cells.Item(122, 9).Value = 1982.3914  # I122
$ws.Cells.Item(122, 10).Value = 2549.6  # J122
$ws.Cells.Item(122, 11).Value = 5947.174199999999  # K122
$ws.Cells.Item(122, 12).Value = 7648.799999999999  # L122
$ws.Cells.Item(122, 13).Value = -3497.174199999999  # M122
$ws.Cells.Item(122, 14).Value = -12548.8  # N122
$ws.Cells.Item(126, 8).Value = 4224.875  # H126
$ws.Cells.Item(126, 9).Value = 3860  # I126
$ws.Cells.Item(126, 10).Value = 4508.6665  # J126
$ws.Cells.Item(126, 11).Value = 11580  # K126
$ws.Cells.Item(126, 12).Value = 13525.9995  # L126
$ws.Cells.Item(126, 13).Value = -9110  # M126
$ws.Cells.Item(126, 14).Value = -18465.9995  # N126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 100565  # H36
$ws.Cells.Item(36, 10).Value = 100565  # J36
$ws.Cells.Item(36, 12).Value = 100565  # L36
$ws.Cells.Item(36, 14).Value = -101689  # N36
$ws.Cells.Item(68, 8).Value = 3499.5  # H68
$ws.Cells.Item(68, 9).Value = 3499.5  # I68
$ws.Cells.Item(68, 10).Value = 0  # J68
$ws.Cells.Item(68, 11).Value = 3499.5  # K68
$ws.Cells.Item(68, 12).Value = 0  # L68
$ws.Cells.Item(68, 13).Value = -2750.5  # M68
$ws.Cells.Item(68, 14).ClearContents()  # N68
$ws.Cells.Item(71, 8).Value = 3499.5  # H71
$ws.Cells.Item(71, 9).Value = 3499.5  # I71
$ws.Cells.Item(71, 10).Value = 0  # J71
$ws.Cells.Item(71, 11).Value = 17497.5  # K71
$ws.Cells.Item(71, 12).Value = 0  # L71
$ws.Cells.Item(71, 13).Value = -13753.5  # M71
$ws.Cells.Item(71, 14).ClearContents()  # N71
$ws.Cells.Item(100, 8).Value = 8897.637000000001  # H100
$ws.Cells.Item(100, 9).Value = 9218.5  # I100
$ws.Cells.Item(100, 11).Value = 9218.5  # K100
$ws.Cells.Item(100, 13).Value = -8677.5  # M100
$ws.Cells.Item(122, 8).Value = 6510  # H122
$ws.Cells.Item(122, 9).Value = 5604.5454  # I122
$ws.Cells.Item(122, 11).Value = 16813.6362  # K122
$ws.Cells.Item(122, 13).Value = -14363.6362  # M122
$ws.Cells.Item(132, 8).Value = 6882428  # H132
$ws.Cells.Item(132, 9).Value = 776988.25  # I132
$ws.Cells.Item(132, 11).Value = 2330964.75  # K132
$ws.Cells.Item(132, 13).Value = -2328434.75  # M132
$ws.Cells.Item(136, 8).Value = 126874.836  # H136
$ws.Cells.Item(136, 9).Value = 81595.766  # I136
$ws.Cells.Item(136, 10).Value = 244600.4  # J136
$ws.Cells.Item(136, 11).Value = 244787.298  # K136
$ws.Cells.Item(136, 12).Value = 733801.2  # L136
$ws.Cells.Item(136, 13).Value = -242237.298  # M136
$ws.Cells.Item(136, 14).Value = -738901.2  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(59, 8).Value = 0  # H59
$ws.Cells.Item(59, 9).Value = 0  # I59
$ws.Cells.Item(59, 11).Value = 0  # K59
$ws.Cells.Item(59, 13).ClearContents()  # M59
$ws.Cells.Item(63, 8).Value = 21613  # H63
$ws.Cells.Item(66, 8).Value = 21613  # H66
$ws.Cells.Item(104, 8).Value = 84646.664  # H104
$ws.Cells.Item(104, 10).Value = 84646.664  # J104
$ws.Cells.Item(104, 12).Value = 84646.664  # L104
$ws.Cells.Item(104, 14).Value = -91634.664  # N104
$ws.Cells.Item(122, 8).Value = 5145.3335  # H122
$ws.Cells.Item(122, 9).Value = 3094.1667  # I122
$ws.Cells.Item(122, 10).Value = 11298.833  # J122
$ws.Cells.Item(122, 11).Value = 9282.500100000001  # K122
$ws.Cells.Item(122, 12).Value = 33896.499  # L122
$ws.Cells.Item(122, 13).Value = -6832.500100000001  # M122
$ws.Cells.Item(122, 14).Value = -38796.499  # N122
$ws.Cells.Item(132, 8).Value = 281322.8  # H132
$ws.Cells.Item(132, 9).Value = 2055.8276  # I132
$ws.Cells.Item(132, 10).Value = 1438286  # J132
$ws.Cells.Item(132, 11).Value = 6167.4828  # K132
$ws.Cells.Item(132, 12).Value = 4314858  # L132
$ws.Cells.Item(132, 13).Value = -3637.4828  # M132
$ws.Cells.Item(132, 14).Value = -4319918  # N132
$ws.Cells.Item(136, 8).Value = 2591.5454  # H136
$ws.Cells.Item(136, 9).Value = 945.2222  # I136
$ws.Cells.Item(136, 11).Value = 2835.6666  # K136
$ws.Cells.Item(136, 13).Value = -285.6666  # M136
